$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 21
$ws.Range("H3").Value = 24

# Row 5
$ws.Range("E5").Value = 38
$ws.Range("F5").Value = 20
$ws.Range("H5").Value = 24

# Row 6
$ws.Range("E6").Value = 66
$ws.Range("F6").Value = 33
$ws.Range("H6").Value = 40

# Row 10
$ws.Range("E10").Value = 38
$ws.Range("F10").Value = 18
$ws.Range("H10").Value = 20

# Row 15
$ws.Range("E15").Value = 127

# Row 16
$ws.Range("E16").Value = 352

# Row 17
$ws.Range("E17").Value = 41
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 17

# Row 18
$ws.Range("E18").Value = 111
